$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# Row 3
$ws.Range("D3").Value = 10.73

# Row 4
$ws.Range("C4").Value = 9.27
$ws.Range("E4").Value = 10.06

# Row 5
$ws.Range("D5").Value = 9.94
$ws.Range("F5").Value = 10.25

# Row 6
$ws.Range("E6").Value = 9.75
$ws.Range("G6").Value = 10.49
$ws.Range("I6").Value = 9.58

# Row 7
$ws.Range("F7").Value = 9.51
$ws.Range("H7").Value = 9.94
$ws.Range("J7").Value = 9

# Row 8
$ws.Range("G8").Value = 10.06

# Row 9
$ws.Range("F9").Value = 10.42

# Row 10
$ws.Range("G10").Value = 11
